$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-obsolete rows 5-7 (data previously describing the FAPs->MuSCs,
# MuSCs->ECs and MuSCs->MuSCs pairs are no longer present after the TPM update)
$ws.Rows("5:7").Delete()

# --- Row 2: Sending cluster ECs, Target cluster now MuSCs ---
$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 4).Value2 = "MuSCs"
$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 0.3747116666666667
$ws.Cells.Item(2, 8).Value2 = 1.124135
$ws.Cells.Item(2, 9).Value2 = 0.3914669751594584
$ws.Cells.Item(2, 10).Value2 = 0.3914669751594584
$ws.Cells.Item(2, 13).Value2 = 0.0002903333333333334
$ws.Cells.Item(2, 14).Value2 = 0.000871
$ws.Cells.Item(2, 15).Value2 = 1
$ws.Cells.Item(2, 16).Value2 = 1
$ws.Cells.Item(2, 17).Value2 = 0.0001087912872222223
$ws.Cells.Item(2, 18).Value2 = 0.0009791215850000002
$ws.Cells.Item(2, 19).Value2 = 0.3914669751594584
$ws.Cells.Item(2, 20).Value2 = 0.3914669751594584

# --- Row 3: Sending cluster FAPs, Target cluster MuSCs ---
$ws.Cells.Item(3, 1).Value2 = "FAPs"
$ws.Cells.Item(3, 4).Value2 = "MuSCs"
$ws.Cells.Item(3, 5).Value2 = 2
$ws.Cells.Item(3, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(3, 7).Value2 = 0.401547
$ws.Cells.Item(3, 8).Value2 = 1.204641
$ws.Cells.Item(3, 9).Value2 = 0.4195022558883632
$ws.Cells.Item(3, 10).Value2 = 0.4195022558883631
$ws.Cells.Item(3, 13).Value2 = 0.0002903333333333334
$ws.Cells.Item(3, 14).Value2 = 0.000871
$ws.Cells.Item(3, 15).Value2 = 1
$ws.Cells.Item(3, 16).Value2 = 1
$ws.Cells.Item(3, 17).Value2 = 0.000116582479
$ws.Cells.Item(3, 18).Value2 = 0.001049242311
$ws.Cells.Item(3, 19).Value2 = 0.4195022558883632
$ws.Cells.Item(3, 20).Value2 = 0.4195022558883631

# --- Row 4: Sending cluster MuSCs, Target cluster MuSCs ---
$ws.Cells.Item(4, 1).Value2 = "MuSCs"
$ws.Cells.Item(4, 4).Value2 = "MuSCs"
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 0.18094
$ws.Cells.Item(4, 8).Value2 = 0.54282
$ws.Cells.Item(4, 9).Value2 = 0.1890307689521785
$ws.Cells.Item(4, 10).Value2 = 0.1890307689521785
$ws.Cells.Item(4, 13).Value2 = 0.0002903333333333334
$ws.Cells.Item(4, 14).Value2 = 0.000871
$ws.Cells.Item(4, 15).Value2 = 1
$ws.Cells.Item(4, 16).Value2 = 1
$ws.Cells.Item(4, 17).Value2 = [double]"5.253291333333334E-05"
$ws.Cells.Item(4, 18).Value2 = 0.00047279622
$ws.Cells.Item(4, 19).Value2 = 0.1890307689521785
$ws.Cells.Item(4, 20).Value2 = 0.1890307689521785
